# Update "Elapsed Duration(Hrs)" values (column G) across the regional
# sheets (R1, R2, R4, R5, R6) of the Active_Outages workbook.
# Values are stored as plain text (e.g. "3927:28:03"), so we force text
# assignment to avoid Excel reinterpreting them as time values.

$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = "R1"; Cell = "G2"; Value = "3927:49:05" },
    @{ Sheet = "R1"; Cell = "G3"; Value = "67:21:43" },
    @{ Sheet = "R2"; Cell = "G2"; Value = "12109:12:45" },
    @{ Sheet = "R2"; Cell = "G3"; Value = "3238:56:14" },
    @{ Sheet = "R2"; Cell = "G4"; Value = "477:07:48" },
    @{ Sheet = "R4"; Cell = "G2"; Value = "2955:02:34" },
    @{ Sheet = "R4"; Cell = "G3"; Value = "182:14:49" },
    @{ Sheet = "R4"; Cell = "G4"; Value = "70:27:14" },
    @{ Sheet = "R4"; Cell = "G5"; Value = "68:04:47" },
    @{ Sheet = "R5"; Cell = "G2"; Value = "429:01:33" },
    @{ Sheet = "R6"; Cell = "G2"; Value = "69:33:51" }
)

foreach ($change in $changes) {
    $ws = $wb.Worksheets.Item($change.Sheet)
    $ws.Range($change.Cell).Value = $change.Value
}
